{"js": "// Commit: \"update 1 vergelijking met clc bio en ncbi blast toegevoegd\"\n//\n// The document has a single empty paragraph (Dutch \"nl-NL\" language,\n// sitting between the \"Verder is het programma...\" paragraph and the\n// \"8 uur programma\" paragraph). This change fills that empty paragraph\n// with a new sentence comparing the NCBI BLAST website defaults with\n// CLC bio's defaults. Several informal/technical terms (blast, fasta,\n// ncbi, blasten, default, clc, bio) are not in the Dutch dictionary, so\n// Word's spell checker wraps each of them in a <w:proofErr> \"spellStart\"\n// / \"spellEnd\" pair, exactly as it did for similar terms elsewhere in\n// the document (e.g. \"Torrent\", \"blasten\", \"plug-ins\", ...).\n\n// The sentence, split into plain-text segments and dictionary-flagged\n// \"word\" segments (in document order) so the run/proofErr structure can\n// be rebuilt programmatically instead of hand-writing 25 near-identical\n// runs.\nconst segments = [\n  [\"text\", \"Omdat mijn programma geen trim en \"],\n  [\"word\", \"blast\"],\n  [\"text\", \" functie heeft ben ik de \"],\n  [\"word\", \"fasta\"],\n  [\"text\", \" files met de \"],\n  [\"word\", \"blast\"],\n  [\"text\", \" op de \"],\n  [\"word\", \"ncbi\"],\n  [\"text\", \" website gaan \"],\n  [\"word\", \"blasten\"],\n  [\"text\", \". De \"],\n  [\"word\", \"ncbi\"],\n  [\"text\", \" website gebruik andere \"],\n  [\"word\", \"default\"],\n  [\"text\", \" instellingen dan \"],\n  [\"word\", \"clc\"],\n  [\"text\", \" \"],\n  [\"word\", \"bio\"],\n  [\"text\", \". De instellingen van \"],\n  [\"word\", \"clc\"],\n  [\"text\", \" \"],\n  [\"word\", \"bio\"],\n  [\"text\", \" hadden als verschil dat de score voor een match 1 is terwijl de \"],\n  [\"word\", \"ncbi\"],\n  [\"text\", \" 2 gebruikt.\"],\n];\n\nfunction xmlEscape(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\")\n    .replace(/\"/g, \"&quot;\");\n}\n\n// Build the <w:r>/<w:proofErr> run sequence for the paragraph body.\nlet runsXml = \"\";\nfor (const [kind, text] of segments) {\n  const preserveAttr = /^\\s|\\s$/.test(text) ? ' xml:space=\"preserve\"' : \"\";\n  const run =\n    \"<w:r><w:rPr><w:lang w:val=\\\"nl-NL\\\"/></w:rPr>\" +\n    \"<w:t\" + preserveAttr + \">\" + xmlEscape(text) + \"</w:t></w:r>\";\n  if (kind === \"word\") {\n    runsXml += '<w:proofErr w:type=\"spellStart\"/>' + run + '<w:proofErr w:type=\"spellEnd\"/>';\n  } else {\n    runsXml += run;\n  }\n}\n\n// Wrap the new paragraph content in a full OOXML package, as required by\n// Paragraph.insertOoxml. Keep the original paragraph's rsid attributes so\n// the <w:p> opening tag itself is left untouched by the edit.\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p w:rsidR=\"00A8567C\" w:rsidRDefault=\"00A8567C\">' +\n  '<w:pPr><w:rPr><w:lang w:val=\"nl-NL\"/></w:rPr></w:pPr>' +\n  runsXml +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\n// Locate the target paragraph: the single empty paragraph in the\n// document (the Dutch section placeholder between \"Verder is het\n// programma...\" and \"8 uur programma\").\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the empty target paragraph\");\n}\n\n// Replace the (empty) paragraph contents with the fully-formed\n// paragraph, preserving its OOXML structure (runs + proofErr markers).\ntarget.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Commit: \"update 1 vergelijking met clc bio en ncbi blast toegevoegd\"\n#\n# The document has a single empty paragraph (Dutch \"nl-NL\" language,\n# sitting between the \"Verder is het programma...\" paragraph and the\n# \"8 uur programma\" paragraph). This change fills that empty paragraph\n# with a new sentence comparing the NCBI BLAST website defaults with\n# CLC bio's defaults. Several informal/technical terms (blast, fasta,\n# ncbi, blasten, default, clc, bio) are not in the Dutch dictionary, so\n# Word's spell checker wraps each of them in a proofErr \"spellStart\" /\n# \"spellEnd\" pair, exactly as it did for similar terms elsewhere in the\n# document (e.g. \"Torrent\", \"blasten\", \"plug-ins\", ...).\n\n$d = $word.ActiveDocument\n\n# The sentence, split into plain-text segments and dictionary-flagged\n# \"word\" segments (in document order), so the run/proofErr structure can\n# be rebuilt programmatically instead of hand-writing 25 near-identical\n# runs.\n$segments = @(\n    ,(\"text\", \"Omdat mijn programma geen trim en \")\n    ,(\"word\", \"blast\")\n    ,(\"text\", \" functie heeft ben ik de \")\n    ,(\"word\", \"fasta\")\n    ,(\"text\", \" files met de \")\n    ,(\"word\", \"blast\")\n    ,(\"text\", \" op de \")\n    ,(\"word\", \"ncbi\")\n    ,(\"text\", \" website gaan \")\n    ,(\"word\", \"blasten\")\n    ,(\"text\", \". De \")\n    ,(\"word\", \"ncbi\")\n    ,(\"text\", \" website gebruik andere \")\n    ,(\"word\", \"default\")\n    ,(\"text\", \" instellingen dan \")\n    ,(\"word\", \"clc\")\n    ,(\"text\", \" \")\n    ,(\"word\", \"bio\")\n    ,(\"text\", \". De instellingen van \")\n    ,(\"word\", \"clc\")\n    ,(\"text\", \" \")\n    ,(\"word\", \"bio\")\n    ,(\"text\", \" hadden als verschil dat de score voor een match 1 is terwijl de \")\n    ,(\"word\", \"ncbi\")\n    ,(\"text\", \" 2 gebruikt.\")\n)\n\nfunction Xml-Escape([string]$s) {\n    $s = $s -replace \"&\", \"&amp;\"\n    $s = $s -replace \"<\", \"&lt;\"\n    $s = $s -replace \">\", \"&gt;\"\n    $s = $s -replace '\"', \"&quot;\"\n    return $s\n}\n\n# Build the <w:r>/<w:proofErr> run sequence for the paragraph body.\n$runsXml = \"\"\nforeach ($seg in $segments) {\n    $kind = $seg[0]\n    $text = $seg[1]\n    $escaped = Xml-Escape $text\n    if ($text -match \"^\\s|\\s$\") {\n        $tOpen = '<w:t xml:space=\"preserve\">'\n    } else {\n        $tOpen = '<w:t>'\n    }\n    $run = '<w:r><w:rPr><w:lang w:val=\"nl-NL\"/></w:rPr>' + $tOpen + $escaped + '</w:t></w:r>'\n    if ($kind -eq \"word\") {\n        $runsXml += '<w:proofErr w:type=\"spellStart\"/>' + $run + '<w:proofErr w:type=\"spellEnd\"/>'\n    } else {\n        $runsXml += $run\n    }\n}\n\n# Wrap the new paragraph content in a full OOXML package, as required by\n# Range.InsertXML. Keep the original paragraph's rsid attributes so the\n# <w:p> opening tag itself is left untouched by the edit.\n$xml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" w:rsidR=\"00A8567C\" w:rsidRDefault=\"00A8567C\"><w:pPr><w:rPr><w:lang w:val=\"nl-NL\"/></w:rPr></w:pPr>' + $runsXml + '</w:p>'\n\n# Locate the target paragraph: the single empty paragraph in the\n# document (the Dutch section placeholder between \"Verder is het\n# programma...\" and \"8 uur programma\").\n$paras = $d.Paragraphs\n$target = $null\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    $t = ($p.Range.Text -replace \"[\\r\\a]+$\", \"\")\n    if ($t.Trim() -eq \"\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not locate the empty target paragraph\"\n}\n\n# Replace the (empty) paragraph contents with the fully-formed\n# paragraph, preserving its OOXML structure (runs + proofErr markers).\n$target.Range.InsertXML($xml)\n"}
